$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B date values (stored as text, hence leading apostrophe to
# preserve the "number stored as text" quote-prefix formatting where needed).
# Introduce the new distinct values in the same order they first appear so
# the shared-string table layout matches.
$ws.Range("B7").Value = "'17/06/2023"
$ws.Range("B12").Value = "'2023/06/17"
$ws.Range("B13").Value = "'18/06/2023"
$ws.Range("B14").Value = "'06/17/2023"
$ws.Range("B1").Value = "'19/07/2023"

$ws.Range("B2").Value = "'19/07/2023"
$ws.Range("B3").Value = "'19/07/2023"
$ws.Range("B4").Value = "'19/07/2023"
$ws.Range("B5").Value = "'19/07/2023"
$ws.Range("B6").Value = "'19/07/2023"

$ws.Range("B8").Value = "'17/06/2023"
$ws.Range("B9").Value = "'17/06/2023"
$ws.Range("B10").Value = "'17/06/2023"
$ws.Range("B11").Value = "'17/06/2023"

$ws.Range("B15").Value = "17/06/2023"

$ws.Range("B16").Value = "'18/06/2023"
$ws.Range("B17").Value = "'18/06/2023"

# Update the active selection to match the saved view state.
[void]$ws.Range("B6").Select()
